# Change the built-in "Block Text" (BlockText) paragraph style to use
# more normal block quotes: indented left and right, keeping the same
# base font/size as the rest of the body text instead of the previous
# smaller / different typeface used only for that style.

$d = $word.ActiveDocument
$style = $d.Styles("Block Text")

# Indent the block quote 24pt (480 twips) on both sides, in addition to
# the already-present "no first-line indent" setting.
$style.ParagraphFormat.LeftIndent = 24
$style.ParagraphFormat.RightIndent = 24
